# Applies the weekly re-sort/update of rows 2-23 (Fruta/Damasco subset) per commit
# "Fruta / hortaliza, semanal". Rows 10 and 11 are unchanged; all other rows in the
# 2..23 range are updated in-place (date, variety, quality, volume, prices, unit,
# origin, $/Kg, Kg/unit) to reflect the new weekly snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D=44161; K="Dina"; L="Primera"; M=300; N=20000; O=20500; P=20250; Q="`$/caja 15 kilos"; R="Región Metropolitana"; S=1350; T=15 },
    @{ Row=3; D=44161; K="Dina"; L="Segunda"; M=100; N=18000; O=18500; P=18250; Q="`$/caja 15 kilos"; R="Región Metropolitana"; S=1217; T=15 },
    @{ Row=4; D=44175; K="Castle Brite"; L="Primera"; M=300; N=21000; O=22000; P=21500; Q="`$/caja 18 kilos"; R="Región Metropolitana"; S=1194; T=18 },
    @{ Row=5; D=44553; K="Modesto"; L="Especial"; M=360; N=23000; O=24000; P=23500; Q="`$/caja 16 kilos"; R="Región Metropolitana"; S=1469; T=16 },
    @{ Row=6; D=44553; K="Modesto"; L="Primera"; M=300; N=21000; O=22000; P=21500; Q="`$/caja 16 kilos"; R="Región Metropolitana"; S=1344; T=16 },
    @{ Row=7; D=44553; K="Modesto"; L="Segunda"; M=240; N=17000; O=18000; P=17500; Q="`$/caja 16 kilos"; R="Región Metropolitana"; S=1094; T=16 },
    @{ Row=8; D=44573; K="Modesto"; L="Especial"; M=300; N=20500; O=21000; P=20750; Q="`$/caja 18 kilos"; R="Región Metropolitana"; S=1153; T=18 },
    @{ Row=9; D=44573; K="Modesto"; L="Primera"; M=400; N=17500; O=18000; P=17750; Q="`$/caja 18 kilos"; R="Región Metropolitana"; S=986; T=18 },
    @{ Row=12; D=44160; K="Castle Brite"; L="Primera"; M=240; N=20500; O=21000; P=20750; Q="`$/caja 15 kilos"; R="Región Metropolitana"; S=1383; T=15 },
    @{ Row=13; D=44559; K="Modesto"; L="Especial"; M=400; N=25000; O=26000; P=25500; Q="`$/caja 18 kilos"; R="Región de O'Higgins"; S=1417; T=18 },
    @{ Row=14; D=44559; K="Modesto"; L="Primera"; M=320; N=22000; O=23000; P=22500; Q="`$/caja 18 kilos"; R="Región de O'Higgins"; S=1250; T=18 },
    @{ Row=15; D=44545; K="Castle Brite"; L="Especial"; M=340; N=22500; O=23000; P=22750; Q="`$/caja 18 kilos"; R="Región de O'Higgins"; S=1264; T=18 },
    @{ Row=16; D=44545; K="Castle Brite"; L="Primera"; M=400; N=20500; O=21000; P=20750; Q="`$/caja 18 kilos"; R="Región de O'Higgins"; S=1153; T=18 },
    @{ Row=17; D=44545; K="Castle Brite"; L="Segunda"; M=300; N=15500; O=16000; P=15750; Q="`$/caja 18 kilos"; R="Región de O'Higgins"; S=875; T=18 },
    @{ Row=18; D=44546; K="Castle Brite"; L="Especial"; M=300; N=22500; O=23000; P=22750; Q="`$/caja 18 kilos"; R="Región Metropolitana"; S=1264; T=18 },
    @{ Row=19; D=44546; K="Castle Brite"; L="Primera"; M=300; N=20500; O=21000; P=20750; Q="`$/caja 18 kilos"; R="Región Metropolitana"; S=1153; T=18 },
    @{ Row=20; D=44566; K="Modesto"; L="Especial"; M=100; N=23000; O=24000; P=23500; Q="`$/caja 18 kilos"; R="Región de O'Higgins"; S=1306; T=18 },
    @{ Row=21; D=44566; K="Modesto"; L="Primera"; M=160; N=21000; O=22000; P=21500; Q="`$/caja 18 kilos"; R="Región de O'Higgins"; S=1194; T=18 },
    @{ Row=22; D=44552; K="Castle Brite"; L="Especial"; M=360; N=20000; O=21000; P=20500; Q="`$/caja 18 kilos"; R="Región Metropolitana"; S=1139; T=18 },
    @{ Row=23; D=44552; K="Castle Brite"; L="Primera"; M=280; N=18000; O=19000; P=18500; Q="`$/caja 18 kilos"; R="Región Metropolitana"; S=1028; T=18 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Cells.Item($r, 4).Value  = $u.D   # D: Fecha
    $ws.Cells.Item($r, 11).Value = $u.K   # K: Variedad
    $ws.Cells.Item($r, 12).Value = $u.L   # L: Calidad
    $ws.Cells.Item($r, 13).Value = $u.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value = $u.N   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $u.O   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $u.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $u.Q   # Q: Unidad de comercializacion
    $ws.Cells.Item($r, 18).Value = $u.R   # R: Origen
    $ws.Cells.Item($r, 19).Value = $u.S   # S: Precio $/Kg
    $ws.Cells.Item($r, 20).Value = $u.T   # T: Kg / unidad
}
